$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9509888291358948
$ws.Range("B1").Value = 1.342366695404053
$ws.Range("C1").Value = 2.244058609008789
$ws.Range("D1").Value = -1
$ws.Range("E1").Value = 1.748462319374084
